$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.940.96"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.984.92"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.48"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.05"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.976.49"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.11"
$ws.Range("E11").Value = "  +6.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.47"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.478.80"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("E17").Value = "  +6.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.86"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.000.04"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "426.65"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.54"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("E22").Value = "  +5.31%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.07"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.39"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.21"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +10.46%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.77"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.63"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.04"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0992"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.991"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.90"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0746"
$ws.Range("E36").Value = "  +9.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.79"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.63"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  +6.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.03"
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.772.06"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0350"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.106"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  +5.14%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.78"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.02"
$ws.Range("E48").Value = "  +22.90%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.99"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.30"
$ws.Range("E51").Value = "  +0.44%  "
